$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C ("Квартира") - shifts old C..H to D..I
$ws.Columns("C").Insert()

# Insert 3 new columns at H..J (right after the old "Показание" column,
# now sitting at G) to make room for tariff #2/#3/#4 - shifts old H..I to K..L
$ws.Range("H1:J1").EntireColumn.Insert()

# ---- Row 1: headers ----
$ws.Range("C1").Value = "Квартира"
$ws.Range("G1").Value = "Показание по тарифу №1"
$ws.Range("H1").Value = "Показание по тарифу №2"
$ws.Range("I1").Value = "Показание по тарифу №3"
$ws.Range("J1").Value = "Показание по тарифу №4"
# (K1 "Житель" and L1 "Источник" shifted here automatically, values retained)

# ---- Row 2: data template (i) ----
$ws.Range("C2").Value = "{d.meter[i].unitName}"
$ws.Range("H2").Value = "{d.meter[i].value2}"
$ws.Range("I2").Value = "{d.meter[i].value3}"
$ws.Range("J2").Value = "{d.meter[i].value4}"

# ---- Row 3: data template (i + 1) ----
$ws.Range("C3").Value = "{d.meter[i + 1].unitName}"
$ws.Range("H3").Value = "{d.meter[i + 1].value2}"
$ws.Range("I3").Value = "{d.meter[i + 1].value3}"
$ws.Range("J3").Value = "{d.meter[i + 1].value4}"
